$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.058.95"
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.27%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9985"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6226"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.15%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07369"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.84%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2924"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.54%  "

# Row 11
$ws.Range("E11").Value = "  -0.13%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07662"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.06%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.830.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.05%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.966"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.51%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6634"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.98%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009031"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.89%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.033"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.050.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.19%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.078.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.31%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "225.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.67%  "

# Row 22
$ws.Range("E22").Value = "  -0.93%  "

# Row 23
$ws.Range("E23").Value = "  -0.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.159"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.000"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.01%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.421"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.50%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1358"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.30%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.80%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.496"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.98%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.056"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.46%  "

# Row 32
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.036"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.19%  "

# Row 33
$ws.Range("E33").Value = "  +0.52%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05251"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.68%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.846"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.52%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.153"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.22%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7334"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.81%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.292.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.42%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01785"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.90%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.748"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.55%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.307"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.62%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9010"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.06%  "

# Row 44
$ws.Range("E44").Value = "  -0.35%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.71%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.975.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.92%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "64.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.83%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5113"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.86%  "

# Row 49
$ws.Range("E49").Value = "  -0.39%  "

# Row 50
$ws.Range("E50").Value = "  -3.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3968"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.76%  "
